# parseRecipe root will now return a parsed recipe
# Add a new "recipe" error row (ID 2) to the Errors sheet / recipe table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errors")

# Copy formatting from the last existing row of the "recipe" table (row 10)
# into the new row (row 11) so fonts / number formats / wrap text match.
$ws.Range("A10:C10").Copy() | Out-Null
$ws.Range("A11:C11").PasteSpecial(-4122) | Out-Null

# Column A holds the numeric-looking error id, stored as text (e.g. "1", "2"),
# matching the other tables (general/axios) on this sheet.
$idCell = $ws.Range("A11")
$idCell.NumberFormat = "@"
$idCell.Value = "2"
$idCell.NumberFormat = """recipe.""@"

$ws.Range("B11").Value = "Could not get recipe from URL"
$ws.Range("C11").Value = "The scraper was unable to get a recipe from the URL provided. Check that the link contains a recipe."

# Match the row height used by the rest of this wrapped-text table.
$ws.Rows.Item(11).RowHeight = 30

# Grow the "recipe" table so the new row becomes part of it.
$lo = $ws.ListObjects.Item("recipe")
$lo.Resize($ws.Range("A9:D11")) | Out-Null

# Update selections / active sheet to match the saved workbook state.
$wsCombined = $wb.Worksheets.Item("Combined")
$wsCombined.Activate() | Out-Null
$wsCombined.Range("D6").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("D11").Select() | Out-Null
